$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '''63.027.26'
$ws.Range('E2').Value = '  +7.89%  '
$ws.Range('D3').Formula = '''3.490.35'
$ws.Range('E3').Value = '  +5.81%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Formula = '''415.72'
$ws.Range('E5').Value = '  +3.90%  '
$ws.Range('D6').Formula = '''127.83'
$ws.Range('E6').Value = '  +16.20%  '
$ws.Range('D7').Formula = '''3.481.80'
$ws.Range('E7').Value = '  +5.60%  '
$ws.Range('D8').Formula = '''0.595'
$ws.Range('E8').Value = '  +1.72%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').Formula = '''0.690'
$ws.Range('E10').Value = '  +9.08%  '
$ws.Range('D11').Formula = '''0.127'
$ws.Range('E11').Value = '  +31.10%  '
$ws.Range('D12').Formula = '''41.82'
$ws.Range('E12').Value = '  +5.08%  '
$ws.Range('D13').Formula = '''0.144'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('D14').Formula = '''4.029.18'
$ws.Range('E14').Value = '  +5.25%  '
$ws.Range('D15').Formula = '''8.69'
$ws.Range('E15').Value = '  +3.33%  '
$ws.Range('D16').Formula = '''20.03'
$ws.Range('E16').Value = '  +4.67%  '
$ws.Range('D17').Formula = '''3.473.22'
$ws.Range('E17').Value = '  +5.70%  '
$ws.Range('D18').Formula = '''62.832.32'
$ws.Range('E18').Value = '  +7.87%  '
$ws.Range('D19').Formula = '''1.05'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('D20').Formula = '''10.85'
$ws.Range('E20').Value = '  -1.23%  '
$ws.Range('D21').Formula = '''0.0000137'
$ws.Range('E21').Value = '  +25.83%  '
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Formula = '''82.10'
$ws.Range('E23').Value = '  +9.75%  '
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Formula = '''315.71'
$ws.Range('E24').Value = '  +5.27%  '
$ws.Range('D25').Formula = '''13.12'
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('D27').Formula = '''31.06'
$ws.Range('E27').Value = '  +9.82%  '
$ws.Range('D28').Formula = '''7.84'
$ws.Range('E28').Value = '  +6.94%  '
$ws.Range('D29').Formula = '''8.08'
$ws.Range('E29').Value = '  +1.82%  '
$ws.Range('D30').Formula = '''0.178'
$ws.Range('E30').Value = '  +4.38%  '
$ws.Range('D31').Formula = '''4.33'
$ws.Range('E31').Value = '  -2.02%  '
$ws.Range('D32').Formula = '''0.119'
$ws.Range('E32').Value = '  +5.18%  '
$ws.Range('D33').Formula = '''2.65'
$ws.Range('E33').Value = '  +23.51%  '
$ws.Range('D34').Formula = '''11.67'
$ws.Range('E34').Value = '  +2.97%  '
$ws.Range('D35').Formula = '''42.91'
$ws.Range('E35').Value = '  +3.88%  '
$ws.Range('E36').Value = '  +0.41%  '
$ws.Range('D37').Formula = '''0.0497'
$ws.Range('E37').Value = '  -3.74%  '
$ws.Range('D38').Formula = '''52.53'
$ws.Range('E38').Value = '  +1.36%  '
$ws.Range('D39').Formula = '''3.55'
$ws.Range('E39').Value = '  +1.59%  '
$ws.Range('D40').Formula = '''0.995'
$ws.Range('E40').Value = '  -0.62%  '
$ws.Range('E41').Value = '  -6.85%  '
$ws.Range('D42').Formula = '''2.03'
$ws.Range('E42').Value = '  +7.57%  '
$ws.Range('E43').Value = '  +2.51%  '
$ws.Range('D44').Formula = '''135.95'
$ws.Range('E44').Value = '  -1.40%  '
$ws.Range('E45').Value = '  +2.14%  '
$ws.Range('D46').Formula = '''17.08'
$ws.Range('E46').Value = '  +1.17%  '
$ws.Range('D47').Formula = '''3.96'
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range('E48').Value = '  -1.75%  '
$ws.Range('D49').Formula = '''22.09'
$ws.Range('E49').Value = '  -2.29%  '
$ws.Range('D50').Formula = '''2.218.58'
$ws.Range('E50').Value = '  +2.21%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Formula = '''3.823.99'
$ws.Range('E51').Value = '  +5.36%  '
